# #5: property aircraft done
# Fix property_category labels that were incorrectly set to "land" for the
# "building" (建物) and "car" (汽車) sheets.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: property_category column I, rows 2-4 -> "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"

# 汽車 (car) sheet: property_category column H, row 2 -> "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
